# ReleaseTracker.xlsx - "Add SV Promo ultra rares"
#
# Rows 18/19 ("Paradox Rift" / set, and "Oinkologne ex Box" / holo promos) got
# their ultra-rare promos added to the tracker, so they flip from the
# "not out" (orange) status to the "waiting for images" (yellow) status, each
# with a short note in column E recording what was added.
#
# Rows 16/17 ("Greninja & Kangaskhan ex Battle Decks" and "Charizard ex
# Premium Collection") had previously been "waiting for images" (yellow);
# their cards are now fully added to the tracker, so they flip to
# "cards added to tracker" (green). Row 17's category note is simplified
# from "holo promos + ultra rare promos" down to just "ultra rare promos"
# now that the plain holo promos are already accounted for elsewhere.
#
# Row 20 ("Gyarados ex Premium Collection") was "not out" (orange) and is now
# confirmed to need "no cards to add" (blue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Carry the current "waiting for images" (yellow) look of row 16/17
#        over to row 18/19 BEFORE row 16/17 themselves switch to green. ---
$ws.Range("A16").Copy()
$ws.Range("A18:A19").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B16:D16").Copy()
$ws.Range("B18:D18").PasteSpecial(-4122)

$ws.Range("B17:D17").Copy()
$ws.Range("B19:D19").PasteSpecial(-4122)

# --- 2. Row 20 moves from "not out" (orange) to "no cards to add" (blue),
#        matching the existing blue look used on row 3. ---
$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial(-4122)

$ws.Range("B3:D3").Copy()
$ws.Range("B20:D20").PasteSpecial(-4122)

# --- 3. Row 16/17 move from "waiting for images" (yellow) to
#        "cards added to tracker" (green), matching row 1's look. ---
$ws.Range("A1").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

$ws.Range("B1:D1").Copy()
$ws.Range("B16:D16").PasteSpecial(-4122)
$ws.Range("B17:D17").PasteSpecial(-4122)

# --- 4. Content edits: category tweaks + new "what was added" notes. ---
$ws.Range("D17").Value = "ultra rare promos"
$ws.Range("D19").Value = "holo promos"
$ws.Range("E18").Value = "ultra rares added"
$ws.Range("E19").Value = "cosmos Lechonk"

# --- 5. Selection / scroll position left where the edits were made. ---
$ws.Range("E19").Select()
